$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts Venue..geometry from B:M to C:N)
$ws.Columns.Item(2).Insert()

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# The inserted column picks up column A's formatting; reset the data rows
# (row 2 down) back to the default/unstyled look used by the other data
# columns, while the header row keeps the bold header style.
$ws.Range($ws.Cells.Item(2, 2), $ws.Cells.Item($lastRow, 2)).Style = "Normal"

# New column B mirrors column A's row index values ("Unnamed: 0" from the
# pandas DataFrame index column), with the same bold/bordered header style
# used by the rest of row 1 (copy formats only from the neighboring header
# cell, which already carries that style).
$ws.Cells.Item(1, 2).Value2 = "Unnamed: 0"
$ws.Cells.Item(1, 3).Copy()
$ws.Cells.Item(1, 2).PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value2 = $ws.Cells.Item($r, 1).Value2
}

$ws.Range("A1").Select()
